# ---------------------------------------------------------------------------
# Capital adequacy 2021 Q4 - add Table_2 worksheet with capital ratio summary
# and drop the now-redundant empty inline-string cells on Table_1.
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Table_1: clear the empty placeholder cells that no longer appear -------
$ws1.Range("B2").Value = ""
$ws1.Range("A3").Value = ""
$ws1.Range("B37").Value = ""

# --- Insert the new "Table_2" sheet right after "Table_1" -------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Table_2"

# --- Header row (bold / centered / boxed, matching Table_1's header style) --
$ws2.Range("A1").Value = "Əmsal"
$ws2.Range("B1").Value = "Norma (Sistem əhəmiyyətli)"
$ws2.Range("C1").Value = "Norma (Banklar istisna)"
$ws2.Range("D1").Value = "Fakt"

$ws1.Range("A1").Copy() | Out-Null
$ws2.Range("A1:D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Data rows ----------------------------------------------------------
# Percentage-looking values must stay plain text, not be auto-converted to
# numbers, so force a text number format before writing them and then
# restore the default ("Normal") style so no stray format sticks around.
$textCells = @("B2", "C2", "D2", "B3", "C3", "D3", "D4")
foreach ($addr in $textCells) {
    $ws2.Range($addr).NumberFormat = "@"
}

$ws2.Range("A2").Value = "9.  I dərəcəli  kapitalın  adekvatlıq əmsalı"
$ws2.Range("B2").Value = "6.0%"
$ws2.Range("C2").Value = "5.0%"
$ws2.Range("D2").Value = "9.5%"

$ws2.Range("A3").Value = "10. məcmu kapitalın  adekvatlıq  əmsalı"
$ws2.Range("B3").Value = "12.0%"
$ws2.Range("C3").Value = "10.0%"
$ws2.Range("D3").Value = "15.5%"

$ws2.Range("A4").Value = "11. Leverec əmsalı"
$ws2.Range("B4").Value = "minimum 5%"
$ws2.Range("C4").Value = "minimum 4%"
$ws2.Range("D4").Value = "5.7%"

foreach ($addr in $textCells) {
    $ws2.Range($addr).Style = "Normal"
}

# --- Leave the original sheet active, as before the edit --------------------
$ws1.Activate()
